# Apply indicator-code renumbering across the quality-metrics workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# SCHEME_MEASURES: MQMSxx -> MQME00x
# ---------------------------------------------------------------------
$wsScheme = $wb.Worksheets.Item("SCHEME_MEASURES")
$wsScheme.Range("A2").Value = "MQME001"
$wsScheme.Range("A3").Value = "MQME002"
$wsScheme.Range("A4").Value = "MQME003"
$wsScheme.Range("A5").Value = "MQME004"
$wsScheme.Range("A6").Value = "MQME005"

# ---------------------------------------------------------------------
# METADATA_ISSUES: MQME11/MQME12/MQME01/MQME14 -> MQME013/014/008/009
# ---------------------------------------------------------------------
$wsIssues = $wb.Worksheets.Item("METADATA_ISSUES")
$wsIssues.Range("A2:A3").Value = "MQME013"
$wsIssues.Range("A4:A156").Value = "MQME014"
$wsIssues.Range("A157:A207").Value = "MQME008"
$wsIssues.Range("A208").Value = "MQME009"

# ---------------------------------------------------------------------
# METADATA_MEASURES: drop the "Total number of columns" row, renumber
# the remaining two indicators.
# ---------------------------------------------------------------------
$wsMeasures = $wb.Worksheets.Item("METADATA_MEASURES")
$wsMeasures.Rows.Item(2).Delete()
$wsMeasures.Range("A2").Value = "MQME006"
$wsMeasures.Range("A3").Value = "MQME007"

# ---------------------------------------------------------------------
# METADATA_METRICS: renumber IQMExx -> MQID0xx, update descriptions /
# values, and add four new indicator rows at the bottom.
# ---------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("METADATA_METRICS")
$defaultStyle = $wsMetrics.Range("A2").Style

function Set-TextValue($cell, $val) {
    # Force the cell to stay a plain text string (rather than Excel
    # auto-converting a "NN.NN%" literal into a formatted percentage
    # number), while keeping the default (unstyled) cell format.
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = $defaultStyle
}

$wsMetrics.Range("A2").Value = "MQID001"
$wsMetrics.Range("B2").Value = "Table names in singular"
Set-TextValue $wsMetrics.Range("C2") "100.00%"

$wsMetrics.Range("A3").Value = "MQID002"
$wsMetrics.Range("B3").Value = "Table with recommended name length"
Set-TextValue $wsMetrics.Range("C3") "98.63%"

$wsMetrics.Range("A4").Value = "MQID003"
$wsMetrics.Range("B4").Value = "Columns with correct prefixes"
Set-TextValue $wsMetrics.Range("C4") "86.82%"

$wsMetrics.Range("A5").Value = "MQID004"
$wsMetrics.Range("B5").Value = "Columns with recommended name size"
Set-TextValue $wsMetrics.Range("C5") "100.00%"

$wsMetrics.Range("A6").Value = "MQID005"
$wsMetrics.Range("B6").Value = "Columns with comments"
Set-TextValue $wsMetrics.Range("C6") "95.61%"

$wsMetrics.Range("A7").Value = "MQID006"
$wsMetrics.Range("B7").Value = "Table with standard PK prefixes"
Set-TextValue $wsMetrics.Range("C7") "99.56%"

$wsMetrics.Range("A8").Value = "MQID007"
$wsMetrics.Range("B8").Value = "Table with standard FK prefixes"
Set-TextValue $wsMetrics.Range("C8") "100.00%"

$wsMetrics.Range("A9").Value = "MQID008"
$wsMetrics.Range("B9").Value = "Table with standard UK prefixes"
Set-TextValue $wsMetrics.Range("C9") "100.00%"

$wsMetrics.Range("A10").Value = "MQID009"
$wsMetrics.Range("B10").Value = "NUMBER columns with valid scale"
Set-TextValue $wsMetrics.Range("C10") "100.00%"

$wsMetrics.Range("A11").Value = "MQID010"
$wsMetrics.Range("B11").Value = "Columns with valid num_distinct"
Set-TextValue $wsMetrics.Range("C11") "100.00%"

$wsMetrics.Range("A12").Value = "MQID011"
$wsMetrics.Range("B12").Value = "Columns with valid num_nulls"
Set-TextValue $wsMetrics.Range("C12") "100.00%"
